$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price column (D) while keeping it text.
# Values such as "596.73" look like genuine numbers to Excel and would
# silently be converted (and "167.20" would lose its trailing zero), so
# cells whose new text parses as a plain float are forced to the "@"
# (text) number format first. Values like "67.678.55" (two dots) are
# never number-parsable and can be assigned directly.
function Set-PriceText {
    param($range, [string]$text)

    # A simple decimal number (optional sign, digits, optional single
    # ".digits") parses as a real number in Excel; values with two dots
    # (e.g. "67.678.55") never do and are safe to assign directly.
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $text
}

# row -> @{ D = '<new price text>'; E = '<new volume text>' }
$updates = @{
    2  = @{ D = '67.678.55'; E = '  +0.46%  ' }
    3  = @{ D = '3.800.87';  E = '  +0.93%  ' }
    4  = @{ E = '  +0.00%  ' }
    5  = @{ D = '596.73';    E = '  +0.66%  ' }
    6  = @{ D = '167.20';    E = '  +0.95%  ' }
    7  = @{ E = '  +0.00%  ' }
    8  = @{ D = '0.519';     E = '  +0.63%  ' }
    9  = @{ E = '  +1.30%  ' }
    10 = @{ D = '6.31';      E = '  -0.87%  ' }
    11 = @{ D = '0.449';     E = '  +0.32%  ' }
    12 = @{ E = '  -0.07%  ' }
    13 = @{ E = '  +0.34%  ' }
    14 = @{ D = '4.442.51';  E = '  +0.99%  ' }
    15 = @{ D = '3.797.05';  E = '  +0.06%  ' }
    16 = @{ D = '18.59';     E = '  +4.90%  ' }
    17 = @{ D = '67.709.85'; E = '  +0.49%  ' }
    18 = @{ E = '  +2.14%  ' }
    19 = @{ E = '  +0.30%  ' }
    20 = @{ D = '461.49';    E = '  +1.27%  ' }
    21 = @{ D = '9.96';      E = '  -2.34%  ' }
    22 = @{ D = '0.701';     E = '  +0.79%  ' }
    23 = @{ E = '  +0.64%  ' }
    24 = @{ D = '83.44' }
    25 = @{ E = '  +2.46%  ' }
    26 = @{ E = '  -0.39%  ' }
    27 = @{ E = '  +0.08%  ' }
    28 = @{ E = '  +0.14%  ' }
    29 = @{ D = '3.943.80';  E = '  +0.76%  ' }
    30 = @{ E = '  +0.08%  ' }
    31 = @{ D = '2.24';      E = '  +3.28%  ' }
    32 = @{ D = '7.31';      E = '  +1.70%  ' }
    33 = @{ D = '29.57';     E = '  -0.14%  ' }
    34 = @{ E = '  -0.02%  ' }
    35 = @{ E = '  -0.72%  ' }
    36 = @{ D = '3.741.58';  E = '  +0.58%  ' }
    37 = @{ D = '0.100';     E = '  +0.72%  ' }
    38 = @{ D = '3.37';      E = '  +2.34%  ' }
    40 = @{ E = '  +0.48%  ' }
    41 = @{ E = '  +0.96%  ' }
    42 = @{ D = '1.00';      E = '  -0.04%  ' }
    44 = @{ E = '  +3.05%  ' }
    45 = @{ E = '  +1.20%  ' }
    46 = @{ D = '43.02';     E = '  -1.46%  ' }
    47 = @{ D = '8.34';      E = '  +0.16%  ' }
    48 = @{ D = '27.37';     E = '  +8.88%  ' }
    49 = @{ E = '  -0.06%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('D')) {
        Set-PriceText $ws.Range("D$row") $vals['D']
    }
    if ($vals.ContainsKey('E')) {
        $ws.Range("E$row").Value = $vals['E']
    }
}

# Rows 50 and 51 swap coin identities (ONDO <-> Bittensor) along with
# their own refreshed price / volume figures.
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-PriceText $ws.Range("D50") '395.90'
$ws.Range("E50").Value = '  +1.95%  '

$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-PriceText $ws.Range("D51") '1.35'
$ws.Range("E51").Value = '  +11.11%  '
